# Regenerate save_data to use K instead of Strike#: update column G (K) values
# for rows 2-9 on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2 = 4
    3 = 3
    4 = 6
    5 = 4
    6 = 3
    7 = 5
    8 = 4
    9 = 2
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
